$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.654690618762475
$ws.Range("C2").Value = 0.696400625978091
$ws.Range("D2").Value = 0.819672131147541
$ws.Range("E2").Value = 0.660377358490566
$ws.Range("F2").Value = 0.643600180913614

$ws.Range("B3").Value = 0.826086956521739
$ws.Range("C3").Value = 0.808970099667774
$ws.Range("D3").Value = 0.846625766871166
$ws.Range("E3").Value = 0.757633587786259
$ws.Range("F3").Value = 0.675473359753413

$ws.Range("B4").Value = 0.760869565217391
$ws.Range("C4").Value = 0.796099290780142
$ws.Range("D4").Value = 0.872689938398357
$ws.Range("E4").Value = 0.648230088495575
$ws.Range("F4").Value = 0.640974639482844
